$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New purchase records (rows 17-19), continuing the shuttlecock_buy log
# with a new brand "RSL Gold".

$rows = @(
    @{ Row = 17; Index = 16; Date = 45528; Brand = "RSL Gold"; PriceRod = 600; NRod = 6; Fee = 0 },
    @{ Row = 18; Index = 17; Date = 45600; Brand = "RSL Gold"; PriceRod = 600; NRod = 2; Fee = 0 },
    @{ Row = 19; Index = 18; Date = 45627; Brand = "RSL Gold"; PriceRod = 600; NRod = 2; Fee = 0 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Index
    $ws.Cells.Item($row, 2).Value = $r.Date
    $ws.Cells.Item($row, 2).NumberFormat = "d-mmm-yy"
    $ws.Cells.Item($row, 3).Value = $r.Brand
    $ws.Cells.Item($row, 4).Value = $r.PriceRod
    $ws.Cells.Item($row, 5).Value = $r.NRod
    $ws.Cells.Item($row, 6).Value = $r.Fee
    $ws.Cells.Item($row, 7).Formula = "=D$row*E$row+F$row"
    $ws.Cells.Item($row, 8).Formula = "=E$row*12"
    $ws.Cells.Item($row, 9).Formula = "=ROUNDUP(G$row/H$row,0)"
}

$ws.Range("D23").Select()
